$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 592091
$ws.Range("C2").Value = -0.004900305531919003
$ws.Range("B3").Value = 587032
$ws.Range("C3").Value = -0.008581390604376793
$ws.Range("B4").Value = 594678
$ws.Range("C4").Value = 0.01294092461466789
$ws.Range("B5").Value = 622732
$ws.Range("C5").Value = 0.04609593376517296
$ws.Range("B6").Value = 636010
$ws.Range("C6").Value = 0.02109854109585285
$ws.Range("B7").Value = 634981
$ws.Range("C7").Value = -0.001619063667021692
$ws.Range("B8").Value = 629714
$ws.Range("C8").Value = -0.008328898809850216
$ws.Range("B9").Value = 627068
$ws.Range("C9").Value = -0.004210436251014471
$ws.Range("B10").Value = 631229
$ws.Range("C10").Value = 0.00661390321329236
$ws.Range("B11").Value = 633555
$ws.Range("C11").Value = 0.003678879700601101
$ws.Range("B12").Value = 632502
$ws.Range("C12").Value = -0.001662912545725703
$ws.Range("B13").Value = 627845
$ws.Range("C13").Value = -0.007390342652797699
$ws.Range("B14").Value = 626907
$ws.Range("C14").Value = -0.001495314529165626
$ws.Range("B15").Value = 625399
$ws.Range("C15").Value = -0.002408578991889954
$ws.Range("B16").Value = 627828
$ws.Range("C16").Value = 0.003876081435009837
$ws.Range("B17").Value = 630155
$ws.Range("C17").Value = 0.003700263798236847
$ws.Range("B18").Value = 632391
$ws.Range("C18").Value = 0.003541411831974983
$ws.Range("B19").Value = 634004
$ws.Range("C19").Value = 0.002547688782215118
$ws.Range("B20").Value = 634045
$ws.Range("C20").Value = 0.00006418088742066175
$ws.Range("B21").Value = 638199
$ws.Range("C21").Value = 0.006530053913593292
$ws.Range("B22").Value = 637741
$ws.Range("C22").Value = -0.0007174259517341852
$ws.Range("B23").Value = 629357
$ws.Range("C23").Value = -0.01323303300887346
$ws.Range("B24").Value = 631768
$ws.Range("C24").Value = 0.003823841223493218
$ws.Range("B25").Value = 629879
$ws.Range("C25").Value = -0.002994842361658812
$ws.Range("B26").Value = 633204
$ws.Range("C26").Value = 0.005265426822006702
$ws.Range("B27").Value = 637285
$ws.Range("C27").Value = 0.006425098981708288
$ws.Range("B28").Value = 633987
$ws.Range("C28").Value = -0.005188418086618185
$ws.Range("B29").Value = 631765
$ws.Range("C29").Value = -0.003510715439915657
$ws.Range("B30").Value = 630361
$ws.Range("C30").Value = -0.002224794821813703
$ws.Range("B31").Value = 634017
$ws.Range("C31").Value = 0.00578342517837882
$ws.Range("B32").Value = 633177
$ws.Range("C32").Value = -0.001326139783486724
$ws.Range("B33").Value = 625531
$ws.Range("C33").Value = -0.01214920729398727
$ws.Range("B34").Value = 622915
$ws.Range("C34").Value = -0.004191330634057522
$ws.Range("B35").Value = 625226
$ws.Range("C35").Value = 0.003702787915244699
$ws.Range("B36").Value = 630674
$ws.Range("C36").Value = 0.008675815537571907
$ws.Range("B37").Value = 631872
$ws.Range("C37").Value = 0.001898121903650463
$ws.Range("B38").Value = 627885
$ws.Range("C38").Value = -0.006329057272523642
$ws.Range("B39").Value = 631581
$ws.Range("C39").Value = 0.005869032815098763
$ws.Range("B40").Value = 635209
$ws.Range("C40").Value = 0.005727432668209076
$ws.Range("B41").Value = 629330
$ws.Range("C41").Value = -0.009297622367739677
$ws.Range("B42").Value = 627493
$ws.Range("C42").Value = -0.00292280176654458
$ws.Range("B43").Value = 629883
$ws.Range("C43").Value = 0.003801406361162663
$ws.Range("B44").Value = 632861
$ws.Range("C44").Value = 0.004716643132269382
$ws.Range("B45").Value = 628497
$ws.Range("C45").Value = -0.006918886676430702
$ws.Range("B46").Value = 629485
$ws.Range("C46").Value = 0.001571498229168355
$ws.Range("B47").Value = 622664
$ws.Range("C47").Value = -0.01089569460600615
$ws.Range("B48").Value = 627286
$ws.Range("C48").Value = 0.00739624397829175
$ws.Range("B49").Value = 623812
$ws.Range("C49").Value = -0.005552992690354586
$ws.Range("B50").Value = 625022
$ws.Range("C50").Value = 0.001938386703841388
$ws.Range("B51").Value = 619836
$ws.Range("C51").Value = -0.008332218043506145
$ws.Range("B52").Value = 620964
$ws.Range("C52").Value = 0.001818535034544766
$ws.Range("B53").Value = 615622
$ws.Range("C53").Value = -0.008640664629638195
$ws.Range("B54").Value = 610330
$ws.Range("C54").Value = -0.008633963763713837
$ws.Range("B55").Value = 610323
$ws.Range("C55").Value = -0.00001112596783059416
$ws.Range("B56").Value = 609382
$ws.Range("C56").Value = -0.001542193465866148
$ws.Range("B57").Value = 613031
$ws.Range("C57").Value = 0.005970211233943701
$ws.Range("B58").Value = 614824
$ws.Range("C58").Value = 0.002919837366789579
$ws.Range("B59").Value = 615854
$ws.Range("C59").Value = 0.001674434752203524
$ws.Range("B60").Value = 615721
$ws.Range("C60").Value = -0.000215426945942454
$ws.Range("B61").Value = 615426
$ws.Range("C61").Value = -0.0004790779494214803
$ws.Range("B62").Value = 618288
$ws.Range("C62").Value = 0.004640447907149792
$ws.Range("B63").Value = 619164
$ws.Range("C63").Value = 0.001415095292031765
$ws.Range("B64").Value = 620238
$ws.Range("C64").Value = 0.001733743702061474
$ws.Range("B65").Value = 617556
$ws.Range("C65").Value = -0.004333783406764269
$ws.Range("B66").Value = 617939
$ws.Range("C66").Value = 0.0006206942489370704
$ws.Range("B67").Value = 620760
$ws.Range("C67").Value = 0.004555430728942156
$ws.Range("B68").Value = 620569
$ws.Range("C68").Value = -0.0003071909595746547
$ws.Range("B69").Value = 615955
$ws.Range("C69").Value = -0.007462609559297562
$ws.Range("B70").Value = 616727
$ws.Range("C70").Value = 0.001253113965503871
$ws.Range("B71").Value = 619192
$ws.Range("C71").Value = 0.003989258781075478
$ws.Range("B72").Value = 619961
$ws.Range("C72").Value = 0.001241911319084466
$ws.Range("B73").Value = 620987
$ws.Range("C73").Value = 0.001653402228839695
$ws.Range("B74").Value = 622089
$ws.Range("C74").Value = 0.001773793832398951
$ws.Range("B75").Value = 621679
$ws.Range("C75").Value = -0.0006596596795134246
$ws.Range("B76").Value = 620120
$ws.Range("C76").Value = -0.002511280123144388
$ws.Range("B77").Value = 621691
$ws.Range("C77").Value = 0.002529565012082458
$ws.Range("B78").Value = 621806
$ws.Range("C78").Value = 0.0001856459712143987
$ws.Range("B79").Value = 625149
$ws.Range("C79").Value = 0.005361692514270544
$ws.Range("B80").Value = 622988
$ws.Range("C80").Value = -0.003462595632299781
$ws.Range("B81").Value = 622560
$ws.Range("C81").Value = -0.0006872529629617929
$ws.Range("B82").Value = 626333
$ws.Range("C82").Value = 0.0060426932759583
$ws.Range("B83").Value = 625523
$ws.Range("C83").Value = -0.001293898443691432
$ws.Range("B84").Value = 623546
$ws.Range("C84").Value = -0.003165753558278084
$ws.Range("B85").Value = 622265
$ws.Range("C85").Value = -0.002056481316685677
$ws.Range("B86").Value = 625856
$ws.Range("C86").Value = 0.005753915291279554
$ws.Range("B87").Value = 627232
$ws.Range("C87").Value = 0.002196959452703595
$ws.Range("B88").Value = 625665
$ws.Range("C88").Value = -0.002501975512132049
$ws.Range("B89").Value = 628317
$ws.Range("C89").Value = 0.004230115562677383
$ws.Range("B90").Value = 627718
$ws.Range("C90").Value = -0.0009538685553707182
$ws.Range("B91").Value = 608960
$ws.Range("C91").Value = -0.03033816628158092
$ws.Range("B92").Value = 610203
$ws.Range("C92").Value = 0.002039549173787236
$ws.Range("B93").Value = 609414
$ws.Range("C93").Value = -0.001293927896767855
$ws.Range("B94").Value = 609430
$ws.Range("C94").Value = 0.00002577547638793476
$ws.Range("B95").Value = 613690
$ws.Range("C95").Value = 0.006965355481952429
$ws.Range("B96").Value = 609774
$ws.Range("C96").Value = -0.006401112768799067
$ws.Range("B97").Value = 608980
$ws.Range("C97").Value = -0.001303470693528652
$ws.Range("B98").Value = 612687
$ws.Range("C98").Value = 0.006069098133593798
$ws.Range("B99").Value = 615352
$ws.Range("C99").Value = 0.004339731764048338
$ws.Range("B100").Value = 609257
$ws.Range("C100").Value = -0.009954815730452538
$ws.Range("B101").Value = 607702
$ws.Range("C101").Value = -0.002555917017161846

Write-Host "Updated rows 2-101 columns B and C"
